# Rename the "Population ID" header (and its linked table column) to
# "Population Name", widen column B to fit the longer text, and move the
# active selection to B2 - matching the author's manual header edit on
# the PopRecordings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updating the header cell's value also renames the linked table column.
$ws.Range("B1").Value = "Population Name"

# Widen column B so the longer header text still fits (was auto best-fit).
$ws.Columns.Item(2).ColumnWidth = 18

# Move the active cell/selection to B2, matching the updated selection.
$ws.Range("B2").Select()
